# Apply updated cryptocurrency price/volume figures scraped on
# Wed Sep  6 20:58:28 UTC 2023 (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.738.09'
$ws.Cells.Item(2, 5).Value = '  -0.27%  '
$ws.Cells.Item(3, 4).Value = '1.632.38'
$ws.Cells.Item(3, 5).Value = '  -0.09%  '
$ws.Cells.Item(4, 5).Value = '  -0.28%  '
$ws.Cells.Item(5, 4).Value = '''215.08'
$ws.Cells.Item(5, 5).Value = '  +0.04%  '
$ws.Cells.Item(6, 4).Value = '''0.503'
$ws.Cells.Item(6, 5).Value = '  -0.84%  '
$ws.Cells.Item(7, 5).Value = '  -0.14%  '
$ws.Cells.Item(8, 5).Value = '  +0.20%  '
$ws.Cells.Item(9, 5).Value = '  -1.03%  '
$ws.Cells.Item(10, 4).Value = '''19.62'
$ws.Cells.Item(10, 5).Value = '  -3.15%  '
$ws.Cells.Item(11, 5).Value = '  +1.14%  '
$ws.Cells.Item(12, 5).Value = '  -0.23%  '
$ws.Cells.Item(13, 4).Value = '1.858.04'
$ws.Cells.Item(13, 5).Value = '  -0.07%  '
$ws.Cells.Item(14, 4).Value = '1.633.33'
$ws.Cells.Item(14, 5).Value = '  -0.32%  '
$ws.Cells.Item(15, 4).Value = '''0.556'
$ws.Cells.Item(15, 5).Value = '  -0.22%  '
$ws.Cells.Item(16, 4).Value = '0.0₃0763'
$ws.Cells.Item(16, 5).Value = '  -0.15%  '
$ws.Cells.Item(17, 4).Value = '''62.62'
$ws.Cells.Item(17, 5).Value = '  -0.89%  '
$ws.Cells.Item(18, 4).Value = '25.747.46'
$ws.Cells.Item(18, 5).Value = '  -0.24%  '
$ws.Cells.Item(19, 5).Value = '  -0.20%  '
$ws.Cells.Item(20, 5).Value = '  +1.82%  '
$ws.Cells.Item(21, 4).Value = '''193.77'
$ws.Cells.Item(21, 5).Value = '  +1.06%  '
$ws.Cells.Item(22, 4).Value = '''9.93'
$ws.Cells.Item(22, 5).Value = '  +0.29%  '
$ws.Cells.Item(23, 5).Value = '  +2.35%  '
$ws.Cells.Item(24, 5).Value = '  -0.16%  '
$ws.Cells.Item(25, 5).Value = '  +3.23%  '
$ws.Cells.Item(26, 4).Value = '''142.31'
$ws.Cells.Item(26, 5).Value = '  +2.39%  '
$ws.Cells.Item(27, 5).Value = '  -0.75%  '
$ws.Cells.Item(28, 4).Value = '''6.86'
$ws.Cells.Item(28, 5).Value = '  +0.74%  '
$ws.Cells.Item(29, 4).Value = '''15.50'
$ws.Cells.Item(29, 5).Value = '  -0.02%  '
$ws.Cells.Item(30, 5).Value = '  -0.08%  '
$ws.Cells.Item(31, 4).Value = '''0.0491'
$ws.Cells.Item(31, 5).Value = '  -0.54%  '
$ws.Cells.Item(32, 5).Value = '  +0.98%  '
$ws.Cells.Item(33, 4).Value = '''3.23'
$ws.Cells.Item(33, 5).Value = '  -0.34%  '
$ws.Cells.Item(34, 5).Value = '  +0.43%  '
$ws.Cells.Item(35, 5).Value = '  -0.01%  '
$ws.Cells.Item(36, 4).Value = '''0.900'
$ws.Cells.Item(36, 5).Value = '  +0.08%  '
$ws.Cells.Item(37, 4).Value = '1.126.74'
$ws.Cells.Item(37, 5).Value = '  -0.31%  '
$ws.Cells.Item(38, 5).Value = '  -1.65%  '
$ws.Cells.Item(39, 4).Value = '''0.545'
$ws.Cells.Item(39, 5).Value = '  -1.71%  '
$ws.Cells.Item(40, 5).Value = '  -0.92%  '
$ws.Cells.Item(41, 5).Value = '  +0.58%  '
$ws.Cells.Item(42, 5).Value = '  +2.07%  '
$ws.Cells.Item(43, 4).Value = '''99.58'
$ws.Cells.Item(43, 5).Value = '  +0.72%  '
$ws.Cells.Item(44, 4).Value = '''0.804'
$ws.Cells.Item(44, 5).Value = '  +0.50%  '
$ws.Cells.Item(45, 4).Value = '1.767.61'
$ws.Cells.Item(46, 5).Value = '  +0.86%  '
$ws.Cells.Item(47, 4).Value = '''54.95'
$ws.Cells.Item(47, 5).Value = '  -1.03%  '
$ws.Cells.Item(48, 5).Value = '  -2.26%  '
$ws.Cells.Item(49, 5).Value = '  +0.11%  '
$ws.Cells.Item(50, 5).Value = '  +3.62%  '
$ws.Cells.Item(51, 4).Value = '''7.58'
$ws.Cells.Item(51, 5).Value = '  -2.65%  '
